$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- Sheet "About": header/version text in A2 and citation text in A6 ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# --- Sheet "Boundaries and methane sources": build_version column (S2:S32) ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 32; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $current = $cell.Value()
    $cell.Value = $current.Replace($oldStamp, $newStamp)
}
